# Update the correct week: shift week labels and refresh seasonality index values
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Map of row -> new week label (column A) and new seasonality index (column P)
$updates = @(
    @{ Row = 2;  Week = "W05"; Seasonality = 0.84 },
    @{ Row = 3;  Week = "W06"; Seasonality = 1 },
    @{ Row = 4;  Week = "W07"; Seasonality = 1.18 },
    @{ Row = 5;  Week = "W08"; Seasonality = $null },
    @{ Row = 6;  Week = "W09"; Seasonality = 1.01 },
    @{ Row = 7;  Week = "W10"; Seasonality = 0.84 },
    @{ Row = 8;  Week = "W11"; Seasonality = 0.85 },
    @{ Row = 9;  Week = "W12"; Seasonality = 1.18 },
    @{ Row = 10; Week = "W13"; Seasonality = 0.8 },
    @{ Row = 11; Week = "W14"; Seasonality = 0.83 },
    @{ Row = 12; Week = "W15"; Seasonality = 1.12 },
    @{ Row = 13; Week = "W16"; Seasonality = 0.98 },
    @{ Row = 14; Week = "W17"; Seasonality = $null },
    @{ Row = 15; Week = "W18"; Seasonality = 0.91 },
    @{ Row = 16; Week = "W19"; Seasonality = 1.05 },
    @{ Row = 17; Week = "W20"; Seasonality = 0.8100000000000001 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.Week
    if ($null -ne $u.Seasonality) {
        $ws.Cells.Item($u.Row, 16).Value = $u.Seasonality
    }
}
